# Update crypto price/volume data (and a few swapped coin rows) per the
# "Updated cryptos list" GitHub Actions commit.
#
# Cells in column D that look like plain numbers ("5.46", "165.20", etc.)
# are forced to stay as text (NumberFormat "@") before assignment so that
# Excel does not silently convert them to numeric values and drop
# significant trailing zeros / formatting, matching the original
# inline-string cell type used throughout this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.101.85"
$ws.Range("E2").Value = "  +2.56%  "

# Row 3
$ws.Range("D3").Value = "2.546.77"
$ws.Range("E3").Value = "  +3.11%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.66"
$ws.Range("E5").Value = "  +1.86%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.67"
$ws.Range("E6").Value = "  +1.81%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.572"
$ws.Range("E8").Value = "  +0.99%  "

# Row 9
$ws.Range("D9").Value = "2.565.54"
$ws.Range("E9").Value = "  +3.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  +2.68%  "

# Row 11
$ws.Range("E11").Value = "  +1.59%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.46"
$ws.Range("E12").Value = "  +0.33%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.362"
$ws.Range("E13").Value = "  +4.54%  "

# Row 14
$ws.Range("D14").Value = "2.996.05"
$ws.Range("E14").Value = "  +3.31%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.19"
$ws.Range("E15").Value = "  +2.98%  "

# Row 16
$ws.Range("D16").Value = "60.056.00"
$ws.Range("E16").Value = "  +2.71%  "

# Row 17
$ws.Range("E17").Value = "  +5.39%  "

# Row 18
$ws.Range("D18").Value = "2.506.34"
$ws.Range("E18").Value = "  +1.53%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.27"
$ws.Range("E19").Value = "  +0.61%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.34"
$ws.Range("E20").Value = "  +2.58%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.13"
$ws.Range("E21").Value = "  +2.28%  "

# Row 22
$ws.Range("E22").Value = "  +0.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.95"
$ws.Range("E23").Value = "  +4.57%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.12"
$ws.Range("E24").Value = "  +4.88%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.436"
$ws.Range("E25").Value = "  +1.53%  "

# Row 26
$ws.Range("E26").Value = "  +5.10%  "

# Row 27
$ws.Range("E27").Value = "  -0.53%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.01"
$ws.Range("E28").Value = "  +4.77%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.13"
$ws.Range("E29").Value = "  +6.79%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0797"
$ws.Range("E30").Value = "  +6.22%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.82"
$ws.Range("E31").Value = "  +3.09%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.19"
$ws.Range("E32").Value = "  -2.03%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "165.20"

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.48"
$ws.Range("E34").Value = "  +6.73%  "

# Row 35
$ws.Range("E35").Value = "  +0.20%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.75"
$ws.Range("E36").Value = "  +1.92%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.44"
$ws.Range("E37").Value = "  +3.58%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.63"
$ws.Range("E38").Value = "  +3.84%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.03"
$ws.Range("E39").Value = "  +1.26%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.60"
$ws.Range("E40").Value = "  -2.00%  "

# Row 41
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "301.88"
$ws.Range("E41").Value = "  +0.62%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.839"
$ws.Range("E42").Value = "  +8.64%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.73"
$ws.Range("E43").Value = "  +3.61%  "

# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.995"
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.609"
$ws.Range("E45").Value = "  +3.25%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.16"
$ws.Range("E47").Value = "  +2.72%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.00"
$ws.Range("E48").Value = "  +3.38%  "

# Row 49
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0938"
$ws.Range("E49").Value = "  +2.07%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0521"
$ws.Range("E50").Value = "  +2.24%  "

# Row 51
$ws.Range("E51").Value = "  +2.38%  "

